{"js": "// Remove \"Standard\" from document titles:\n//   \"Square One Standard Offer Letter\"  -> \"Square One Offer Letter\"\n//   \"Square One Standard Employee Terms\" -> \"Square One Employee Terms\"\nconst body = context.document.body;\n\nconst replacements = [\n  [\"Square One Standard Offer Letter\", \"Square One Offer Letter\"],\n  [\"Square One Standard Employee Terms\", \"Square One Employee Terms\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Remove \"Standard\" from document titles:\n#   \"Square One Standard Offer Letter\"   -> \"Square One Offer Letter\"\n#   \"Square One Standard Employee Terms\" -> \"Square One Employee Terms\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"Square One Standard Offer Letter\"\n$find.Replacement.Text = \"Square One Offer Letter\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Wrap = 1\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, $null, $null, $null, $null, $null, $null, $null, [ref]$find.Replacement.Text, 2)\n\n$find.Text = \"Square One Standard Employee Terms\"\n$find.Replacement.Text = \"Square One Employee Terms\"\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, $null, $null, $null, $null, $null, $null, $null, [ref]$find.Replacement.Text, 2)\n"}
